$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4793.2354
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 4780.3125
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 4780.3125
$ws.Range("M40").Value = -4825
$ws.Range("N40").Value = -5130.3125

# Row 61
$ws.Range("H61").Value = 924.75
$ws.Range("I61").Value = 924.75
$ws.Range("K61").Value = 2774.25
$ws.Range("M61").Value = -2602.25

# Row 132
$ws.Range("H132").Value = 3264.3125
$ws.Range("I132").Value = 3269.1
$ws.Range("K132").Value = 9807.299999999999
$ws.Range("M132").Value = -7277.299999999999

# Row 134
$ws.Range("H134").Value = 80000
$ws.Range("I134").Value = 80000
$ws.Range("K134").Value = 80000
$ws.Range("M134").Value = -74930


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2661.5557
$ws.Range("I2").Value = 1803
$ws.Range("J2").Value = 3090.8333
$ws.Range("K2").Value = 1803
$ws.Range("L2").Value = 3090.8333
$ws.Range("M2").Value = -1690
$ws.Range("N2").Value = -3316.8333

# Row 10
$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10340

# Row 13
$ws.Range("H13").Value = 25000000
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 32
$ws.Range("H32").Value = 3469.2952
$ws.Range("I32").Value = 2909.5173
$ws.Range("K32").Value = 2909.5173
$ws.Range("M32").Value = -2622.5173

# Row 116
$ws.Range("H116").Value = 2661.5557
$ws.Range("I116").Value = 1803
$ws.Range("J116").Value = 3090.8333
$ws.Range("K116").Value = 1803
$ws.Range("L116").Value = 3090.8333
$ws.Range("M116").Value = 491
$ws.Range("N116").Value = -7678.8333


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2661.5557
$ws.Range("I3").Value = 1803
$ws.Range("J3").Value = 3090.8333
$ws.Range("K3").Value = 1803
$ws.Range("L3").Value = 3090.8333
$ws.Range("M3").Value = -1689
$ws.Range("N3").Value = -3318.8333

# Row 5
$ws.Range("H5").Value = 1030.8
$ws.Range("I5").Value = 663.5
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 663.5
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = -550.5
$ws.Range("N5").Value = -2726

# Row 94
$ws.Range("H94").Value = 4250
$ws.Range("I94").Value = 3750
$ws.Range("J94").Value = 5750
$ws.Range("K94").Value = 3750
$ws.Range("L94").Value = 5750
$ws.Range("M94").Value = -3299
$ws.Range("N94").Value = -6652

# Row 134
$ws.Range("H134").Value = 3898.9
$ws.Range("I134").Value = 4123.625
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 12370.875
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -9835.875
$ws.Range("N134").Value = -14070


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5735.2
$ws.Range("I31").Value = 1880.375
$ws.Range("J31").Value = 10140.714
$ws.Range("K31").Value = 1880.375
$ws.Range("L31").Value = 10140.714
$ws.Range("M31").Value = -1585.375
$ws.Range("N31").Value = -10730.714

# Row 34
$ws.Range("H34").Value = 5735.2
$ws.Range("I34").Value = 1880.375
$ws.Range("J34").Value = 10140.714
$ws.Range("K34").Value = 1880.375
$ws.Range("L34").Value = 10140.714
$ws.Range("M34").Value = -1678.375
$ws.Range("N34").Value = -10544.714

# Row 86
$ws.Range("H86").Value = 5999.5
$ws.Range("J86").Value = 5999.5
$ws.Range("L86").Value = 5999.5
$ws.Range("N86").Value = -8245.5

# Row 89
$ws.Range("H89").Value = 5999.5
$ws.Range("J89").Value = 5999.5
$ws.Range("L89").Value = 29997.5
$ws.Range("N89").Value = -41229.5

# Row 134
$ws.Range("H134").Value = 4737
$ws.Range("I134").Value = 3934.5
$ws.Range("J134").Value = 8749.5
$ws.Range("K134").Value = 11803.5
$ws.Range("L134").Value = 26248.5
$ws.Range("M134").Value = -9268.5
$ws.Range("N134").Value = -31318.5

# Row 135
$ws.Range("H135").Value = 247152.23
$ws.Range("J135").Value = 247152.23
$ws.Range("L135").Value = 247152.23
$ws.Range("N135").Value = -257292.23


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 6000
$ws.Range("I5").Value = 4000
$ws.Range("J5").Value = 8000
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = -11888
$ws.Range("N5").Value = -24224

# Row 14
$ws.Range("H14").Value = 833
$ws.Range("I14").Value = 833
$ws.Range("K14").Value = 2499
$ws.Range("M14").Value = -2326

# Row 98
$ws.Range("H98").Value = 551.5
$ws.Range("J98").Value = 1000
$ws.Range("L98").Value = 3000
$ws.Range("N98").Value = -5996

# Row 104
$ws.Range("H104").Value = 36666.668
$ws.Range("I104").Value = 30000
$ws.Range("J104").Value = 40000
$ws.Range("K104").Value = 90000
$ws.Range("L104").Value = 120000
$ws.Range("M104").Value = -87379
$ws.Range("N104").Value = -125242

# Row 134
$ws.Range("H134").Value = 3409.577
$ws.Range("I134").Value = 1509.9333
$ws.Range("K134").Value = 4529.7999
$ws.Range("M134").Value = 540.2001

# Row 135
$ws.Range("H135").Value = 6000
$ws.Range("I135").Value = 4000
$ws.Range("J135").Value = 8000
$ws.Range("K135").Value = 36000
$ws.Range("L135").Value = 72000
$ws.Range("M135").Value = -33465
$ws.Range("N135").Value = -77070


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 5094111.5
$ws.Range("I11").Value = 6790479
$ws.Range("K11").Value = 6790479
$ws.Range("M11").Value = -6790340

# Row 102
$ws.Range("H102").Value = 12998.6
$ws.Range("I102").Value = 8748.5
$ws.Range("K102").Value = 8748.5
$ws.Range("M102").Value = -7126.5

# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7059.8
$ws.Range("I7").Value = 5900
$ws.Range("J7").Value = 8799.5
$ws.Range("K7").Value = 5900
$ws.Range("L7").Value = 8799.5
$ws.Range("M7").Value = -5788
$ws.Range("N7").Value = -9023.5

# Row 22
$ws.Range("H22").Value = 3841.48
$ws.Range("I22").Value = 3806.0908
$ws.Range("J22").Value = 3869.2856
$ws.Range("K22").Value = 3806.0908
$ws.Range("L22").Value = 3869.2856
$ws.Range("M22").Value = -3511.0908
$ws.Range("N22").Value = -4459.2856

# Row 27
$ws.Range("H27").Value = 3841.48
$ws.Range("I27").Value = 3806.0908
$ws.Range("J27").Value = 3869.2856
$ws.Range("K27").Value = 3806.0908
$ws.Range("L27").Value = 3869.2856
$ws.Range("M27").Value = -3699.0908
$ws.Range("N27").Value = -4083.2856

# Row 68
$ws.Range("H68").Value = 5294.28
$ws.Range("I68").Value = 2959.3333
$ws.Range("J68").Value = 7449.615
$ws.Range("K68").Value = 2959.3333
$ws.Range("L68").Value = 7449.615
$ws.Range("M68").Value = -2210.3333
$ws.Range("N68").Value = -8947.615

# Row 71
$ws.Range("H71").Value = 5294.28
$ws.Range("I71").Value = 2959.3333
$ws.Range("J71").Value = 7449.615
$ws.Range("K71").Value = 14796.6665
$ws.Range("L71").Value = 37248.075
$ws.Range("M71").Value = -11052.6665
$ws.Range("N71").Value = -44736.075

# Row 93
$ws.Range("H93").Value = 5912.0713
$ws.Range("I93").Value = 2553.8
$ws.Range("K93").Value = 2553.8
$ws.Range("M93").Value = -1305.8

# Row 126
$ws.Range("H126").Value = 7059.8
$ws.Range("I126").Value = 5900
$ws.Range("J126").Value = 8799.5
$ws.Range("K126").Value = 17700
$ws.Range("L126").Value = 26398.5
$ws.Range("M126").Value = -15230
$ws.Range("N126").Value = -31338.5

# Row 132
$ws.Range("H132").Value = 3647.7856
$ws.Range("I132").Value = 3437.0417
$ws.Range("K132").Value = 10311.1251
$ws.Range("M132").Value = -7781.125100000001


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 105
$ws.Range("I8").Value = 105
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 105
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 35
$ws.Range("N8").ClearContents()

# Row 15
$ws.Range("H15").Value = 33333.332
$ws.Range("J15").Value = 33333.332
$ws.Range("L15").Value = 33333.332
$ws.Range("N15").Value = -33909.332

# Row 113
$ws.Range("H113").Value = 1682.25
$ws.Range("I113").Value = 1312.625
$ws.Range("K113").Value = 3937.875
$ws.Range("M113").Value = -1767.875

# Row 122
$ws.Range("H122").Value = 4197.8335
$ws.Range("I122").Value = 2558.6365
$ws.Range("J122").Value = 5584.846
$ws.Range("K122").Value = 7675.9095
$ws.Range("L122").Value = 16754.538
$ws.Range("M122").Value = -5225.9095
$ws.Range("N122").Value = -21654.538

# Row 126
$ws.Range("H126").Value = 1686.8235
$ws.Range("I126").Value = 1248.1538
$ws.Range("K126").Value = 3744.4614
$ws.Range("M126").Value = -1274.4614

# Row 132
$ws.Range("H132").Value = 1561.6316
$ws.Range("I132").Value = 1570.44
$ws.Range("J132").Value = 1544.6923
$ws.Range("K132").Value = 4711.32
$ws.Range("L132").Value = 4634.0769
$ws.Range("M132").Value = -2181.32
$ws.Range("N132").Value = -9694.0769

# Row 136
$ws.Range("H136").Value = 1134.2858
$ws.Range("I136").Value = 1161.1538
$ws.Range("J136").Value = 785
$ws.Range("K136").Value = 3483.4614
$ws.Range("L136").Value = 2355
$ws.Range("M136").Value = -933.4614000000001
$ws.Range("N136").Value = -7455

